$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new "City, Country" column right after each Team Member's
# "Affiliation Email" column (processed left-to-right so each insertion
# point already accounts for columns shifted by earlier insertions).
$ws.Columns("E:E").Insert()
$ws.Columns("I:I").Insert()
$ws.Columns("M:M").Insert()
$ws.Columns("Q:Q").Insert()
$ws.Columns("U:U").Insert()
$ws.Columns("Y:Y").Insert()

# Fill in the new header cells.
$ws.Range("E1").Value = "City, Country"
$ws.Range("I1").Value = "City, Country"
$ws.Range("M1").Value = "City, Country"
$ws.Range("Q1").Value = "City, Country"
$ws.Range("U1").Value = "City, Country"
$ws.Range("Y1").Value = "City, Country"

# The inserted columns for members 4 and 5 (Q1, U1) inherited the
# neighbouring "explicit black" font style instead of the "theme" font
# style used by the rest of that header row; fix that up to match the
# formatting used by the other member blocks.
$ws.Range("C1").Copy()
$ws.Range("Q1").PasteSpecial(-4122)
$ws.Range("C1").Copy()
$ws.Range("U1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Update the note text: it used to mention city/country explicitly; that
# information now lives in its own column, so trim the note.
$ws.Range("B14").Value = "Affiliation should include department, university /institute/company name in full"

# Move the active selection (cosmetic, matches the saved view state).
$ws.Range("F18").Select()
